$d = $word.ActiveDocument

# Locate the paragraph that currently holds the "_GoBack" bookmark and the
# text "Ler livro do Pavel " (it is the second-to-last paragraph in the
# document, right before the trailing empty paragraph / sectPr).
$bm = $d.Bookmarks("_GoBack")
$para = $bm.Range.Paragraphs.Item(1)
$paraStart = $para.Range.Start
$paraEnd = $para.Range.End
$target = $d.Range($paraStart, $paraEnd)

# Replace that single paragraph with two paragraphs:
#   1) "Ler livro do Pavel " - same formatting as before, bookmark removed
#   2) "Como o artigo de resiliência serve?" - strike only (no rFonts),
#      wrapped by the relocated "_GoBack" bookmark.
$newParaXml = (
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" ' +
            'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
                '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                    '<w:body>' +
                        '<w:p>' +
                            '<w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:strike/></w:rPr></w:pPr>' +
                            '<w:r>' +
                                '<w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:strike/></w:rPr>' +
                                '<w:t xml:space="preserve">Ler livro do Pavel </w:t>' +
                            '</w:r>' +
                        '</w:p>' +
                        '<w:p>' +
                            '<w:pPr><w:rPr><w:strike/></w:rPr></w:pPr>' +
                            '<w:r>' +
                                '<w:rPr><w:strike/></w:rPr>' +
                                '<w:t>Como o artigo de resiliência serve?</w:t>' +
                            '</w:r>' +
                        '</w:p>' +
                    '</w:body>' +
                '</w:document>' +
            '</pkg:xmlData>' +
        '</pkg:part>' +
    '</pkg:package>'
)

$target.InsertXML($newParaXml)

# The paragraph split above dropped the "_GoBack" bookmark (InsertXML
# replaces the whole range it is anchored in), so recreate it around the
# freshly inserted second paragraph's text.
$count = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($count - 1)
$bmStart = $newPara.Range.Start
$bmEnd = $newPara.Range.End - 1
$bmRange = $d.Range($bmStart, $bmEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
